$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ 2 = 0.3397006374899547; 3 = 0.02866471113603097; 5 = 0.3378702794578814; 6 = 2.693714048937238; 7 = 1.20957717133804; 8 = 1.161705211693757; 10 = 0.08968341031793159; 11 = 0.2911386485056653; 13 = 0.3466445029133709 }
    3 = @{ 2 = 0.3101547195475689; 3 = 0.02538525436314387; 5 = 0.3320661613269991; 6 = 2.672381942019086; 7 = 1.205136168216157; 8 = 1.164462089716523; 10 = 0.09034952885836134; 11 = 0.2609142533154483; 13 = 0.3312153550935051 }
    4 = @{ 2 = 0.2921502670799327; 3 = 0.02336004607641939; 5 = 0.328663194587584; 6 = 2.660655456467822; 7 = 1.203073496367537; 8 = 1.166630220449647; 10 = 0.0907837669669469; 11 = 0.2424394326356207; 13 = 0.321922896572147 }
    5 = @{ 2 = 0.2848478912238193; 3 = 0.02253184242987061; 5 = 0.327316954772563; 6 = 2.656221471229003; 7 = 1.202399568869382; 8 = 1.167633227193846; 10 = 0.09096706453278358; 11 = 0.2349318280763839; 13 = 0.3181817895064682 }
    6 = @{ 2 = 0.2836374312105079; 3 = 0.02239414393670813; 5 = 0.3270958605542518; 6 = 2.65550602105705; 7 = 1.202297716796735; 8 = 1.167806989819226; 10 = 0.09099788386665963; 11 = 0.2336864708499746; 13 = 0.3175633415706969 }
    7 = @{ 2 = 0.2920516444256691; 3 = 0.02334888842869276; 5 = 0.3286448746681714; 6 = 2.660594263011944; 7 = 1.203063733371863; 8 = 1.166643263690034; 10 = 0.09078621330927383; 11 = 0.2423380970914337; 13 = 0.3218722577354782 }
    8 = @{ 2 = 0.3294848908588222; 3 = 0.02753636525009995; 5 = 0.3358356851951285; 6 = 2.686073930978679; 7 = 1.207907847749865; 8 = 1.162557073031081; 10 = 0.0899078470919612; 11 = 0.2807001092480448; 13 = 0.3412870026165251 }
    9 = @{ 2 = 0.403973990207902; 3 = 0.03565610803025265; 5 = 0.3512108647640488; 6 = 2.74693802762755; 7 = 1.222695263491104; 8 = 1.158319927038249; 10 = 0.08838581089879582; 11 = 0.356584925250786; 13 = 0.3807940139720642 }
    10 = @{ 2 = 0.4593628245719401; 3 = 0.04156657014928555; 5 = 0.3632830754192469; 6 = 2.798330040017618; 7 = 1.236811842602563; 8 = 1.157515197788513; 10 = 0.0873899828996727; 11 = 0.4127412329370372; 13 = 0.410695171720846 }
    11 = @{ 2 = 0.4847051178155937; 3 = 0.04424366299922156; 5 = 0.3689435531495988; 6 = 2.823166336825068; 7 = 1.243946334978517; 8 = 1.157651810459043; 10 = 0.08696358340786681; 11 = 0.4383769752825799; 13 = 0.4244884718638247 }
    12 = @{ 2 = 0.4943224613570862; 3 = 0.04525574653781916; 5 = 0.3711112635881619; 6 = 2.832781249427285; 7 = 1.246750935127494; 8 = 1.157775930799673; 10 = 0.08680594839057143; 11 = 0.4480974613596231; 13 = 0.4297390765306943 }
    13 = @{ 2 = 0.4922502734797831; 3 = 0.04503785073109157; 5 = 0.3706433325132963; 6 = 2.830701163521923; 7 = 1.246142330031134; 8 = 1.157745978148967; 10 = 0.08683972731497924; 11 = 0.4460034143952498; 13 = 0.4286070486601545 }
    14 = @{ 2 = 0.4854959272449832; 3 = 0.04432696128223768; 5 = 0.3691214071099864; 6 = 2.823953152030526; 7 = 1.244175005965872; 8 = 1.157660570562115; 10 = 0.08695053781604223; 11 = 0.4391764300948466; 13 = 0.424919894077064 }
    15 = @{ 2 = 0.4813613944972417; 3 = 0.04389130291576748; 5 = 0.3681923350500611; 6 = 2.819847150692283; 7 = 1.242983379339861; 8 = 1.15761768605654; 10 = 0.08701891185640953; 11 = 0.4349963659525997; 13 = 0.4226649683300536 }
    16 = @{ 2 = 0.4577095616182021; 3 = 0.04139138189671598; 5 = 0.3629165409271522; 6 = 2.796736291335279; 7 = 1.236359966185034; 8 = 1.157516393916183; 10 = 0.08741838512529654; 11 = 0.41106767423571; 13 = 0.4097975841566495 }
    17 = @{ 2 = 0.4432371032320077; 3 = 0.03985478537764209; 5 = 0.3597231958938494; 6 = 2.782932131263664; 7 = 1.232479581370768; 8 = 1.157583077437252; 10 = 0.08767026889784368; 11 = 0.396411129892897; 13 = 0.4019527371842457 }
    18 = @{ 2 = 0.4349266460371837; 3 = 0.03896988352420294; 5 = 0.3579023569999649; 6 = 2.775129540731712; 7 = 1.230314752430004; 8 = 1.157668741248273; 10 = 0.08781764979273632; 11 = 0.3879895490246668; 13 = 0.3974585830801516 }
    19 = @{ 2 = 0.432115229503097; 3 = 0.03867008352392531; 5 = 0.3572885829351549; 6 = 2.772511271880248; 7 = 1.229593284228756; 8 = 1.157705867849785; 10 = 0.08786798028323517; 11 = 0.3851396067615269; 13 = 0.3959400333762488 }
    20 = @{ 2 = 0.4447763026422251; 3 = 0.0400184717710772; 5 = 0.3600614889295741; 6 = 2.784387405994011; 7 = 1.232885710823055; 8 = 1.157571081793364; 10 = 0.0876431962337314; 11 = 0.3979704679372844; 13 = 0.4027859729223309 }
    21 = @{ 2 = 0.4874792802272054; 3 = 0.04453581219522107; 5 = 0.3695677768893546; 6 = 2.825929506185503; 7 = 1.244750060057584; 8 = 1.157683691430663; 10 = 0.08691788599692174; 11 = 0.4411813350099294; 13 = 0.4260021584275222 }
    22 = @{ 2 = 0.5155091122106512; 3 = 0.04747839960924694; 5 = 0.3759217759752005; 6 = 2.854303544134041; 7 = 1.253104207394614; 8 = 1.158179269384107; 10 = 0.0864662000841534; 11 = 0.4694965994254687; 13 = 0.4413348315493266 }
    23 = @{ 2 = 0.5005380518084053; 3 = 0.04590877931690329; 5 = 0.3725176362607883; 6 = 2.839047696934074; 7 = 1.24859039151616; 8 = 1.157876124897456; 10 = 0.08670522631532052; 11 = 0.4543774446101452; 13 = 0.4331369288251068 }
    24 = @{ 2 = 0.4440803995792351; 3 = 0.03994447380088673; 5 = 0.3599084997413158; 6 = 2.783729060129758; 7 = 1.232701893961988; 8 = 1.157576357602835; 10 = 0.08765542778145807; 11 = 0.3972654766785979; 13 = 0.4024092173475324 }
    25 = @{ 2 = 0.3837066503964763; 3 = 0.03346924819319952; 5 = 0.3469151462522149; 6 = 2.729302706527221; 7 = 1.21812580073663; 8 = 1.159061292469914; 10 = 0.08877608140019078; 11 = 0.3359854165477998; 13 = 0.3699526572446246 }
}

foreach ($row in $data.Keys) {
    foreach ($col in $data[$row].Keys) {
        $ws.Cells.Item($row, $col).Value = $data[$row][$col]
    }
}

Write-Output "applied $($data.Count) rows"